$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the previous table contents (A1:H50); styles (bold/border on headers & labels) are preserved on cells that already carried them ---
$ws.Range("A1:H50").ClearContents()

# --- Header row (B1:I1): d=1 .. d=5, d=6 (new), d=7, d=10 ---
$ws.Range("B1").Value = "d=1"
$ws.Range("C1").Value = "d=2"
$ws.Range("D1").Value = "d=3"
$ws.Range("E1").Value = "d=4"
$ws.Range("F1").Value = "d=5"
$ws.Range("G1").Value = "d=6"
$ws.Range("H1").Value = "d=7"
$ws.Range("I1").Value = "d=10"

# --- Row labels, column A (A2:A57) ---
$ws.Range("A2").Value = "ARMA_I(0,1,0)"
$ws.Range("A3").Value = "ARMA_I(0,1,1)"
$ws.Range("A4").Value = "ARMA_I(0,1,2)"
$ws.Range("A5").Value = "ARMA_I(0,10,0)"
$ws.Range("A6").Value = "ARMA_I(0,10,1)"
$ws.Range("A7").Value = "ARMA_I(0,10,2)"
$ws.Range("A8").Value = "ARMA_I(0,2,0)"
$ws.Range("A9").Value = "ARMA_I(0,2,1)"
$ws.Range("A10").Value = "ARMA_I(0,2,2)"
$ws.Range("A11").Value = "ARMA_I(0,3,0)"
$ws.Range("A12").Value = "ARMA_I(0,3,1)"
$ws.Range("A13").Value = "ARMA_I(0,3,2)"
$ws.Range("A14").Value = "ARMA_I(0,4,0)"
$ws.Range("A15").Value = "ARMA_I(0,4,1)"
$ws.Range("A16").Value = "ARMA_I(0,4,2)"
$ws.Range("A17").Value = "ARMA_I(0,5,0)"
$ws.Range("A18").Value = "ARMA_I(0,5,1)"
$ws.Range("A19").Value = "ARMA_I(0,5,2)"
$ws.Range("A20").Value = "ARMA_I(0,6,0)"
$ws.Range("A21").Value = "ARMA_I(0,6,1)"
$ws.Range("A22").Value = "ARMA_I(0,6,2)"
$ws.Range("A23").Value = "ARMA_I(0,7,0)"
$ws.Range("A24").Value = "ARMA_I(0,7,1)"
$ws.Range("A25").Value = "ARMA_I(0,7,2)"
$ws.Range("A26").Value = "ARMA_I(1,1,0)"
$ws.Range("A27").Value = "ARMA_I(1,1,1)"
$ws.Range("A28").Value = "ARMA_I(1,10,0)"
$ws.Range("A29").Value = "ARMA_I(1,10,1)"
$ws.Range("A30").Value = "ARMA_I(1,2,0)"
$ws.Range("A31").Value = "ARMA_I(1,2,1)"
$ws.Range("A32").Value = "ARMA_I(1,3,0)"
$ws.Range("A33").Value = "ARMA_I(1,3,1)"
$ws.Range("A34").Value = "ARMA_I(1,4,0)"
$ws.Range("A35").Value = "ARMA_I(1,4,1)"
$ws.Range("A36").Value = "ARMA_I(1,5,0)"
$ws.Range("A37").Value = "ARMA_I(1,5,1)"
$ws.Range("A38").Value = "ARMA_I(1,6,0)"
$ws.Range("A39").Value = "ARMA_I(1,6,1)"
$ws.Range("A40").Value = "ARMA_I(1,7,0)"
$ws.Range("A41").Value = "ARMA_I(1,7,1)"
$ws.Range("A42").Value = "ARMA_I(2,1,0)"
$ws.Range("A43").Value = "ARMA_I(2,1,2)"
$ws.Range("A44").Value = "ARMA_I(2,10,0)"
$ws.Range("A45").Value = "ARMA_I(2,10,2)"
$ws.Range("A46").Value = "ARMA_I(2,2,0)"
$ws.Range("A47").Value = "ARMA_I(2,2,2)"
$ws.Range("A48").Value = "ARMA_I(2,3,0)"
$ws.Range("A49").Value = "ARMA_I(2,3,2)"
$ws.Range("A50").Value = "ARMA_I(2,4,0)"
$ws.Range("A51").Value = "ARMA_I(2,4,2)"
$ws.Range("A52").Value = "ARMA_I(2,5,0)"
$ws.Range("A53").Value = "ARMA_I(2,5,2)"
$ws.Range("A54").Value = "ARMA_I(2,6,0)"
$ws.Range("A55").Value = "ARMA_I(2,6,2)"
$ws.Range("A56").Value = "ARMA_I(2,7,0)"
$ws.Range("A57").Value = "ARMA_I(2,7,2)"

# --- Data values ---
$ws.Range("B2").Value = 30.31806102698224
$ws.Range("B3").Value = 30.76153967020114
$ws.Range("B4").Value = 40.25973713831916
$ws.Range("I5").Value = 96.24236789370403
$ws.Range("I6").Value = 96.13098304617974
$ws.Range("I7").Value = 96.22756446181195
$ws.Range("C8").Value = 50.63498013815493
$ws.Range("C9").Value = 50.43470392334577
$ws.Range("C10").Value = 50.70436496210364
$ws.Range("D11").Value = 98.69172919731389
$ws.Range("D12").Value = 98.78537891235521
$ws.Range("D13").Value = 98.97802536907948
$ws.Range("E14").Value = 98.93984262359832
$ws.Range("E15").Value = 98.9978291233006
$ws.Range("E16").Value = 98.93757540209221
$ws.Range("F17").Value = 98.4813239258053
$ws.Range("F18").Value = 98.58154259493311
$ws.Range("F19").Value = 98.44626910233967
$ws.Range("G20").Value = 98.14905295172701
$ws.Range("G21").Value = 98.20988268362727
$ws.Range("G22").Value = 98.04867160915441
$ws.Range("H23").Value = 97.56959798991242
$ws.Range("H24").Value = 97.61725417075434
$ws.Range("H25").Value = 97.61390091135038
$ws.Range("B26").Value = 40.98592402049587
$ws.Range("B27").Value = 49.03583991382323
$ws.Range("I28").Value = 96.30820249441508
$ws.Range("I29").Value = 96.23109626382599
$ws.Range("C30").Value = 50.19379561441237
$ws.Range("C31").Value = 50.30764186347617
$ws.Range("D32").Value = 98.8470680953949
$ws.Range("D33").Value = 98.79434438846096
$ws.Range("E34").Value = 98.95227071897385
$ws.Range("E35").Value = 98.95128926113065
$ws.Range("F36").Value = 98.45407170166315
$ws.Range("F37").Value = 98.58918864167788
$ws.Range("G38").Value = 98.03385597590334
$ws.Range("G39").Value = 98.04776379758653
$ws.Range("H40").Value = 97.50497962783544
$ws.Range("H41").Value = 97.62698278886836
$ws.Range("B42").Value = 36.33344537134072
$ws.Range("B43").Value = 50.80760093129585
$ws.Range("I44").Value = 96.28893688435859
$ws.Range("I45").Value = 96.09449641807879
$ws.Range("C46").Value = 50.94024582960805
$ws.Range("C47").Value = 50.10465028580613
$ws.Range("D48").Value = 98.9599121273158
$ws.Range("D49").Value = 98.79257303335514
$ws.Range("E50").Value = 98.98878201414485
$ws.Range("E51").Value = 98.85244615087099
$ws.Range("F52").Value = 98.57621715884011
$ws.Range("F53").Value = 98.44329350150969
$ws.Range("G54").Value = 98.04256948786806
$ws.Range("G55").Value = 98.01296453597725
$ws.Range("H56").Value = 97.49120598448928
$ws.Range("H57").Value = 97.63882725993936

# --- Re-apply the header/label style (bold font, thin border, centered) to the newly added cells: I1 and A51:A57 ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A51:A57").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A1").Select() | Out-Null
